# Add the new "Haunted Carriage" and "Demon Gates" timer columns so the
# event timers sort properly alongside the existing "Ancient Nightmare"
# columns (H:J).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (row 1), columns K:P
$ws.Range("K1").Value = "Haunted_Carriage_12PM"
$ws.Range("L1").Value = "Haunted_Carriage_830PM"
$ws.Range("M1").Value = "Haunted_Carriage_10PM"
$ws.Range("N1").Value = "Demon_Gates_12PM"
$ws.Range("O1").Value = "Demon_Gates_830PM"
$ws.Range("P1").Value = "Demon_Gates_10PM"

# Match the column widths that Excel's "best fit" produced for the new
# header text (same sizing pattern already used by columns H:J).
$ws.Columns.Item(11).ColumnWidth = 22.33
$ws.Columns.Item(12).ColumnWidth = 23.33
$ws.Columns.Item(13).ColumnWidth = 22.33
$ws.Columns.Item(14).ColumnWidth = 22.33
$ws.Columns.Item(15).ColumnWidth = 23.33
$ws.Columns.Item(16).ColumnWidth = 22.33

# Reflect the new selection left behind on the sheet after the edit.
$ws.Activate() | Out-Null
$ws.Range("N1:P1").Select() | Out-Null
